$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 262, shifting existing rows 262:342 down to 263:343.
$ws.Rows.Item(262).Insert()

# Populate the newly inserted row 262 with the new record.
$ws.Cells.Item(262, 1).Value = 8
$ws.Cells.Item(262, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(262, 3).Value = "Coquimbo"
$ws.Cells.Item(262, 4).Value = 44900
$ws.Cells.Item(262, 5).Value = 4
$ws.Cells.Item(262, 6).Value = 100112012
$ws.Cells.Item(262, 7).Value = "Espinaca"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 800
$ws.Cells.Item(262, 11).Value = 500
$ws.Cells.Item(262, 12).Value = 600
$ws.Cells.Item(262, 13).Value = 550
$ws.Cells.Item(262, 14).Value = '$/atado 300 a 500 gramos'
$ws.Cells.Item(262, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(262, 16).Value = 1100
$ws.Cells.Item(262, 17).Value = 0.5
$ws.Cells.Item(262, 18).Value = "Hortaliza"
